# Junction_Flooding_157.xlsx edit
#   - "custom accuracy" : round the measurement values on row 5 to 2 decimal
#     places (they were stored with 3 decimals before).
#   - "데이터 1000개"     : trim the data set by removing the last data row
#     (row 6), shrinking the used range from A1:AH6 to A1:AH5.
#   - Column J (10th column) is narrowed from width 8 to width 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Round row 5's numeric measurements to 2 decimal places -------------
# (cells that already only needed <=2 decimals - G5, P5, T5, AD5 - are left
#  untouched since their value does not change)
$ws.Range("B5").Value  = 18.74
$ws.Range("C5").Value  = 13.74
$ws.Range("D5").Value  = 1.18
$ws.Range("E5").Value  = 40.73
$ws.Range("F5").Value  = 33.12
$ws.Range("H5").Value  = 53.2
$ws.Range("I5").Value  = 22.69
$ws.Range("J5").Value  = 10.01
$ws.Range("K5").Value  = 14.81
$ws.Range("L5").Value  = 16.34
$ws.Range("M5").Value  = 17.21
$ws.Range("N5").Value  = 4.71
$ws.Range("O5").Value  = 14.66
$ws.Range("Q5").Value  = 12.43
$ws.Range("R5").Value  = 0.81
$ws.Range("S5").Value  = 0.79
$ws.Range("U5").Value  = 40.88
$ws.Range("V5").Value  = 13.53
$ws.Range("W5").Value  = 27.4
$ws.Range("X5").Value  = 14.38
$ws.Range("Y5").Value  = 2.26
$ws.Range("Z5").Value  = 26.29
$ws.Range("AA5").Value = 11.95
$ws.Range("AB5").Value = 10.62
$ws.Range("AC5").Value = 12.5
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 48.12
$ws.Range("AG5").Value = 7.58

# --- 2. Drop the last data row (row 6), shifting nothing else -------------
$ws.Rows.Item(6).Delete()

# --- 3. Narrow column J (index 10) from width 8 to width 7 -----------------
# ColumnWidth is expressed in "characters", which is offset from the stored
# OOXML column width by a constant ~0.83 padding factor on this sheet's font,
# so request 7 - 0.83 = 6.17 to land on a stored width of 7.
$ws.Columns.Item(10).ColumnWidth = 6.17
